$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate() | Out-Null

# --- Fix typo "Assignemnt 5 - Part 2" -> "Assignment 5 - Part 2" ---
$ws.Range("E37").Value = "Assignment 5 - Part 2"

# --- Row 28 (Assignment 3 exercise): mark COMPLETE and set completion date ---
$ws.Range("A28").Value = "COMPLETE"
$ws.Range("B25").Copy() | Out-Null
$ws.Range("B28").PasteSpecial(-4122) | Out-Null
$ws.Range("B28").Value = 44194

# --- Row 37 (Assignment 5 - Part 2 exercise): mark COMPLETE and set completion date ---
$ws.Range("A37").Value = "COMPLETE"
$ws.Range("B34").Copy() | Out-Null
$ws.Range("B37").PasteSpecial(-4122) | Out-Null
$ws.Range("B37").Value = 44194

# --- New rows for the Python Performance Tuning lesson + Assignment 6 exercise ---
$ws.Range("A38").Value = "NOT STARTED"
$ws.Range("C38").Value = 6
$ws.Range("D38").Value = "Software Carpentry Lesson"
$ws.Range("E38").Value = "Python Performance Tuning - Introduction"

$ws.Range("A39").Value = "NOT STARTED"
$ws.Range("C39").Value = 6
$ws.Range("D39").Value = "Software Carpentry Lesson"
$ws.Range("E39").Value = "Python Performance Tuning - cProfile"

$ws.Range("A40").Value = "NOT STARTED"
$ws.Range("C40").Value = 6
$ws.Range("D40").Value = "Software Carpentry Lesson"
$ws.Range("E40").Value = "Python Performance Tuning - Line_Profiler"

$ws.Range("A41").Value = "NOT STARTED"
$ws.Range("C41").Value = 6
$ws.Range("D41").Value = "Exercise"
$ws.Range("E41").Value = "Assignment 6"

# --- Update the view's saved scroll position / selection ---
$ws.Range("C42").Select() | Out-Null
